# Kosten für Gehäuseherstellung ergänzt
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlLeft    = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignLeft
$xlCenterV = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignCenter

# ---------------------------------------------------------------------------
# 1. Header row (row 1) becomes a 4-column header:
#    Bauteil | Preis in € | Menge in kg | Preis pro kg
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 1).Value = "Bauteil"
$ws.Cells.Item(1, 2).Value = "Preis in €"
$ws.Cells.Item(1, 3).Value = "Menge in kg"
$ws.Cells.Item(1, 4).Value = "Preis pro kg"

$ws.Cells.Item(1, 1).Font.Bold = $true
$ws.Cells.Item(1, 1).VerticalAlignment = $xlCenterV

for ($c = 2; $c -le 4; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = $xlLeft
    $cell.VerticalAlignment = $xlCenterV
}

# ---------------------------------------------------------------------------
# 2. Re-order the existing parts list (rows 2-14) to match the new ordering
# ---------------------------------------------------------------------------
$labels = @(
  "Arduino Mega",
  "Powerbank",
  "RFID Reader Set",
  "20 RFID Chips",
  "Micro SD Karte",
  "Micro SD Reader",
  "1-Kanal 5V Relay/Relais Modul 10A - 250VAC",
  "Drucktaster",
  "4x4 Keypad",
  "LCD 16x2 I2C (vorverlötet)",
  "100 Dioden",
  "LED",
  "Pauschale für Kabel und weiteres Zubehör"
)
$prices = @(11.99, 14.24, 4.99, 6.29, 2.99, 3.99, 3.61, 5.95, 6.99, 5.85, 4.9000000000000004, 5.34, 15)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $labels[$i]
    $bcell = $ws.Cells.Item($r, 2)
    $bcell.Value = $prices[$i]
    $bcell.NumberFormat = "0.00"
    $bcell.HorizontalAlignment = $xlLeft
}

# ---------------------------------------------------------------------------
# 3. Blank separator row 15 (kept empty, like before)
# ---------------------------------------------------------------------------
$ws.Cells.Item(15, 1).Value = $null

# ---------------------------------------------------------------------------
# 4. New "Gehäuseherstellung" (3D-printed enclosure) cost rows 16-20
# ---------------------------------------------------------------------------
$ws.Cells.Item(16, 1).Value = "90° Adapter"
$b16 = $ws.Cells.Item(16, 2)
$b16.Value = 4.3899999999999997
$b16.NumberFormat = "0.00"
$b16.HorizontalAlignment = $xlLeft

$ws.Cells.Item(17, 1).Value = "USB Buchse"
$b17 = $ws.Cells.Item(17, 2)
$b17.Value = 6.3
$b17.NumberFormat = "0.00"
$b17.HorizontalAlignment = $xlLeft

$ws.Cells.Item(18, 1).Value = "PLA Schwarz"
$c18 = $ws.Cells.Item(18, 3)
$c18.Value = 0.312
$c18.HorizontalAlignment = $xlLeft
$c18.VerticalAlignment = $xlCenterV
$d18 = $ws.Cells.Item(18, 4)
$d18.Value = 26.99
$d18.HorizontalAlignment = $xlLeft
$d18.VerticalAlignment = $xlCenterV
$b18 = $ws.Cells.Item(18, 2)
$b18.Formula = "=D18*C18"
$b18.NumberFormat = "0.00"
$b18.HorizontalAlignment = $xlLeft

$ws.Cells.Item(19, 1).Value = "PLA Transparent"
$c19 = $ws.Cells.Item(19, 3)
$c19.Value = 0.41399999999999998
$c19.HorizontalAlignment = $xlLeft
$c19.VerticalAlignment = $xlCenterV
$d19 = $ws.Cells.Item(19, 4)
$d19.Value = 21.99
$d19.HorizontalAlignment = $xlLeft
$d19.VerticalAlignment = $xlCenterV
$b19 = $ws.Cells.Item(19, 2)
$b19.Formula = "=D19*C19"
$b19.NumberFormat = "0.00"
$b19.HorizontalAlignment = $xlLeft

$ws.Cells.Item(20, 1).Value = "Pauschale f. Gewinde etc."
$b20 = $ws.Cells.Item(20, 2)
$b20.Value = 5
$b20.NumberFormat = "0.00"
$b20.HorizontalAlignment = $xlLeft

# ---------------------------------------------------------------------------
# 5. Total row (row 22); row 21 stays an empty separator
# ---------------------------------------------------------------------------
$a22 = $ws.Cells.Item(22, 1)
$a22.Value = "Gesamtkosten"
$a22.Font.Bold = $true

$b22 = $ws.Cells.Item(22, 2)
$b22.Formula = "=SUM(B2:B21)"
$b22.NumberFormat = "0.00"
$b22.Font.Bold = $true
$b22.HorizontalAlignment = $xlLeft

# ---------------------------------------------------------------------------
# 6. Column widths for B (resized) and the new C/D columns
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 10.44140625
$ws.Columns.Item(3).ColumnWidth = 11.33203125
$ws.Columns.Item(4).ColumnWidth = 12.6640625

$ws.Range("A1").Select()
